# Grading_Guide.xlsx update — "Added treeview - needs more work."
#
# The checklist on Sheet1 gets a round of progress updates: several
# "Attempted?" checkboxes flip from FALSE to TRUE (their "Potential Score"
# formulas recalculate automatically), and the two "...lowest cost?" sub-rows
# (which are being dropped from the checklist) have their text/checkbox
# cleared out in place. Finally the active selection moves to C27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Mark additional checklist items as attempted -------------------------
$attemptedNow = @("C5", "C6", "C14", "C15", "C16", "C17", "C24", "C25", "C26", "C27")
foreach ($addr in $attemptedNow) {
    $ws.Range($addr).Value = $true
}

# --- Remove the two "...lowest cost?" sub-items ----------------------------
# Row 18: "Waypoint Support in Dijkstras lowest cost?"
$ws.Range("A18").ClearContents()
$ws.Range("C18").ClearContents()

# Row 23: "Avoiding Specified Stations in Dijkstras lowest cost?"
$ws.Range("A23").ClearContents()

# --- Update selection -------------------------------------------------------
$ws.Range("C27").Select()

$excel.Calculate()
